$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Duplicate the "MM1" milestone block (rows 2-8) into new rows
#    20-26 so that column B/C/D values and number formats/styles are
#    preserved exactly (re-using the existing style indexes).
# ------------------------------------------------------------------
$src = $ws.Range("A2:D8")
$dst = $ws.Range("A20")
$src.Copy($dst)

# ------------------------------------------------------------------
# 2. Re-label the freshly duplicated rows as the "MM2" milestone.
# ------------------------------------------------------------------
$ws.Range("A20").Value = "Approval MM2"
$ws.Range("A21").Value = "Approval MM2 LoD"
$ws.Range("A22").Value = "Approval MM2 Version No"
$ws.Range("A23").Value = "Approval MM2 Original Baseline"
$ws.Range("A24").Value = "Approval MM2 Forecast - Actual"
$ws.Range("A25").Value = "Approval MM2 Status"
$ws.Range("A26").Value = "Approval MM2 Notes"

# ------------------------------------------------------------------
# 3. Rename the original "Assurance MM1" labels to "Approval MM1".
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Approval MM1"
$ws.Range("A3").Value = "Approval MM1 LoD"
$ws.Range("A4").Value = "Approval MM1 Version No"
$ws.Range("A5").Value = "Approval MM1 Original Baseline"
$ws.Range("A6").Value = "Approval MM1 Forecast - Actual"
$ws.Range("A7").Value = "Approval MM1 Status"
$ws.Range("A8").Value = "Approval MM1 Notes"

# ------------------------------------------------------------------
# 4. Restore the selected / active cell used by the workbook view.
# ------------------------------------------------------------------
$null = $ws.Range("B38").Select()
